$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix E26: convert from text "544028" to a real number
$ws.Range("E26").Value = 544028

# Add new row 27 with data
$ws.Range("A27").Value = "21/06/2024 07:44:37"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "TATATECH"
$ws.Range("D27").Value = "Tata Technologies Ltd"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "544028"
$ws.Range("E27").ClearFormats()
$ws.Range("F27").Value = -1.02
$ws.Range("G27").Value = 1000
$ws.Range("H27").Value = 2080736
